# Daily attendance processing - 2026-01-28 14:03:37
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet
# lists the users who recorded attendance for a session. For sessions
# recorded by both the automated System process and a human grader,
# the order of the names needs to be flipped from
#   "System, dnasr281@gmail.com"
# to
#   "dnasr281@gmail.com, System"
# This script scans the used range of column G and swaps the order of
# the two names wherever that exact combination is found, leaving every
# other cell (including rows that only say "System" or only the email)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $current = $cell.Value2
    if ($current -eq $oldText) {
        $cell.Value = $newText
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed 'Recorded By' cell(s) in column G."
